$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Melbourne -> Taichung lane, container upsized 20GP -> 40GP, rates rebased
$ws.Range("B4").Value = "Taichung"
$ws.Range("C4").Value = "40GP"

# Numeric-looking values must stay text (as the rest of the sheet does),
# so force text formatting before writing them.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "500"

$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "400"

$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "400"

$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "400"

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "400"

# J4/K4 (Collect / 14 Days) and E4 (400) are unchanged - leave as-is.

# Drop row 5 (Melbourne/Shanghai 40REHC line) entirely.
$ws.Rows("5").Delete()
